$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (GMHO:0000199 / intervention mechanism of action): the previously
# empty "text" placeholder cells E2:R2 and T2:U2 become blank/number-typed
# empty cells (their string content is cleared).
$ws.Range("E2:R2").ClearContents()
$ws.Range("T2:U2").ClearContents()

# Copy row 2's cell formatting/style down onto the new row 3 before filling
# it in, so row 3 picks up the same fill style (s="2") as row 2.
$ws.Range("A2:V2").Copy()
$ws.Range("A3:V3").PasteSpecial(-4122)  # xlPasteFormats

# Add the new row 3 entry: BFO:0000015 / process
$ws.Range("A3").Value = "BFO:0000015"
$ws.Range("B3").Value = "process"
$ws.Range("C3").Value = "p is a process = Def. p is an occurrent that has temporal proper parts and for some time t, p s-depends_on some material entity at t. (axiom label in BFO2 Reference: [083-003])"
$ws.Range("D3").Value = "occurrent"
$ws.Range("S3").Value = "Proposed"
